$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 4 ("Stuart, A. (2005, May 28) ...") edits ---
$para4 = $tr.Paragraphs(4, 1)

# Run 1: update the citation text to "2005a"
$run1 = $para4.Runs(1, 1)
$run1.Text = "Stuart, A. (2005a, May 28). Sudoku solver by Andrew Stuart. "

# Run 2: the hyperlink run -> underline it
$run2 = $para4.Runs(2, 1)
$run2.Font.Underline = -1

# Run 3: trailing space run -> remove (becomes endParaRPr)
$run3 = $para4.Runs(3, 1)
$run3.Text = ""

Write-Host "Para4 now:" $para4.Text

# --- Insert two new paragraphs after paragraph 4 ---
$para4 = $tr.Paragraphs(4, 1)
$inserted = $para4.InsertAfter("Stuart, A. (2005b, June 9). Naked candidates. SudokuWiki.org - Naked Candidates. https://www.sudokuwiki.org/Naked_Candidates#NP `rStuart, A. (2008, April 9). Hidden candidates. SudokuWiki.org - Hidden Candidates. https://www.sudokuwiki.org/Hidden_Candidates#HP ")
Write-Host "Inserted text:" $inserted.Text

